$d = $word.ActiveDocument

# Locate the paragraph that hosts the " m:'doc.html'.fromHTMLURI() " field
# (the only paragraph whose Range contains a Word field).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph containing the field"
}

$rng = $d.Range($target.Range.Start, $target.Range.End - 1)

# Rebuild the paragraph: the fldChar begin/end + the leading/trailing
# instrText space runs are replaced by literal '{' and '}' text runs,
# while every other instrText run becomes an equivalent w:t run carrying
# the same text. The bookmark stays exactly where it was.
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"' +
    ' w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' +
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m</w:t></w:r>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '<w:r><w:t>''</w:t></w:r>' +
    '<w:r><w:t>doc.html</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>''.fromHTMLURI()</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
    '</w:p>'

[void]$rng.InsertXML($xmlFrag)
